# Refresh the cryptocurrency price/volume snapshot (GitHub Actions data pull).
# Only the D (Price) and E (Volume(1h)) columns for rows 2-51 change; everything
# else (coin name, link, row index) is left untouched.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '47.400.15'
$ws.Range("E2").Value = '  +2.65%  '
$ws.Range("D3").Value = '2.503.44'
$ws.Range("E3").Value = '  +2.10%  '
$ws.Range("D5").Value = '''324.99'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.21%  '
$ws.Range("D6").Value = '''110.27'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +5.15%  '
$ws.Range("E7").Value = '  +1.20%  '
$ws.Range("D8").Value = '''1.00'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.00%  '
$ws.Range("E9").Value = '  +0.35%  '
$ws.Range("E10").Value = '  +10.01%  '
$ws.Range("E11").Value = '  +1.17%  '
$ws.Range("E12").Value = '  +1.03%  '
$ws.Range("D13").Value = '''18.50'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +1.30%  '
$ws.Range("E14").Value = '  +1.95%  '
$ws.Range("D15").Value = '2.894.70'
$ws.Range("E15").Value = '  +2.23%  '
$ws.Range("D16").Value = '2.503.18'
$ws.Range("E16").Value = '  +1.61%  '
$ws.Range("D17").Value = '''0.860'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +1.90%  '
$ws.Range("D18").Value = '47.314.95'
$ws.Range("E18").Value = '  +2.84%  '
$ws.Range("D19").Value = '''12.89'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +2.43%  '
$ws.Range("E20").Value = '  +3.99%  '
$ws.Range("E21").Value = '  +0.96%  '
$ws.Range("D22").Value = '''2.70'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +13.43%  '
$ws.Range("D23").Value = '''70.54'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.95%  '
$ws.Range("D24").Value = '''248.31'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.42%  '
$ws.Range("D25").Value = '''2.60'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +3.67%  '
$ws.Range("E26").Value = '  +0.70%  '
$ws.Range("E27").Value = '  -0.12%  '
$ws.Range("D28").Value = '''2.30'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +1.58%  '
$ws.Range("E29").Value = '  +3.98%  '
$ws.Range("D30").Value = '''35.35'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +4.67%  '
$ws.Range("D31").Value = '''0.138'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +9.04%  '
$ws.Range("E32").Value = '  +1.15%  '
$ws.Range("E33").Value = '  -0.16%  '
$ws.Range("D34").Value = '''5.45'
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").Value = '''0.0796'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +4.76%  '
$ws.Range("E36").Value = '  +0.26%  '
$ws.Range("E37").Value = '  +5.27%  '
$ws.Range("E38").Value = '  +3.74%  '
$ws.Range("E39").Value = '  +1.93%  '
$ws.Range("E40").Value = '  +1.53%  '
$ws.Range("D41").Value = '''121.80'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -3.64%  '
$ws.Range("E42").Value = '  -0.60%  '
$ws.Range("D43").Value = '''21.16'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +1.20%  '
$ws.Range("E44").Value = '  +2.50%  '
$ws.Range("D45").Value = '2.001.75'
$ws.Range("E45").Value = '  +1.71%  '
$ws.Range("D46").Value = '''3.10'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +4.46%  '
$ws.Range("E47").Value = '  -0.31%  '
$ws.Range("E48").Value = '  -3.78%  '
$ws.Range("E49").Value = '  -0.59%  '
$ws.Range("E50").Value = '  +4.36%  '
$ws.Range("D51").Value = '''56.85'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +4.08%  '
